$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45766 = 2025-04-19) for every
# data row (2-43). The workbook was refreshed a day later, so bump each of
# these date values forward by one day (-> serial 45767 = 2025-04-20).
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45766) {
        $cell.Value2 = 45767
    }
}
